$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 208
$ws1.Range("F3").Value = 117
$ws1.Range("F5").Value = 994
$ws1.Range("F6").Value = 5525
$ws1.Range("F7").Value = 494
$ws1.Range("F8").Value = 694
$ws1.Range("F9").Value = 953
$ws1.Range("F10").Value = 824
$ws1.Range("F11").Value = 79
$ws1.Range("F17").Value = 1855
$ws1.Range("F18").Value = 1475
$ws1.Range("F19").Value = 923
$ws1.Range("F21").Value = 195
$ws1.Range("F22").Value = 336
$ws1.Range("F23").Value = 555
$ws1.Range("F24").Value = 155
$ws1.Range("F25").Value = 1055
$ws1.Range("F28").Value = 2934
$ws1.Range("F32").Value = 122
$ws1.Range("F34").Value = 386
$ws1.Range("F39").Value = 292
$ws1.Range("F40").Value = 722
$ws1.Range("F41").Value = 91
$ws1.Range("F44").Value = 69

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 193

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 208
$ws4.Range("F4").Value = 117
$ws4.Range("F5").Value = 994
$ws4.Range("F7").Value = 5525
$ws4.Range("F8").Value = 494
$ws4.Range("F9").Value = 694
$ws4.Range("F11").Value = 193
$ws4.Range("F12").Value = 953
$ws4.Range("F13").Value = 824
$ws4.Range("F16").Value = 79
$ws4.Range("F23").Value = 1855
$ws4.Range("F24").Value = 1475
$ws4.Range("F25").Value = 923
$ws4.Range("F26").Value = 195
$ws4.Range("F27").Value = 336
$ws4.Range("F29").Value = 555
$ws4.Range("F30").Value = 155
$ws4.Range("F31").Value = 1055
$ws4.Range("F32").Value = 2934
$ws4.Range("F36").Value = 122
$ws4.Range("F38").Value = 386
$ws4.Range("F42").Value = 292
$ws4.Range("F43").Value = 722
$ws4.Range("F44").Value = 91
$ws4.Range("F46").Value = 69
